$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# F13: add value 1 (100%), with percentage number format "0%" (border already thin box from prior style)
$ws.Range("F13").Value = 1
$ws.Range("F13").NumberFormat = "0%"

# F14: add value 0.9 (90%), keep existing border (left/right thin, bottom double) but add percent format
$ws.Range("F14").Value = 0.9
$ws.Range("F14").NumberFormat = "0%"

# G14: add note text about lookup issue
$ws.Range("G14").Value = "Một số chuỗi HoTen ko tra cứu được ??. Một số chức năng sẽ bổ sung sau."

# Update the active selection to F14 (matches end-of-edit cursor position)
$ws.Range("F14").Select()

$wb.Save()
